# updated charts in kruskal
$wb = $excel.ActiveWorkbook

# --- Update existing CHOR1 / CHOR2 sheets with revised stats ---
$chor1 = $wb.Worksheets.Item("CHOR1")
$chor1.Range("D2").Value = 5.36
$chor1.Range("I2").Value = 0.86
$chor1.Range("D3").Value = 5.77
$chor1.Range("I3").Value = 0.29
$chor1.Range("D4").Value = 4.2

$chor2 = $wb.Worksheets.Item("CHOR2")
$chor2.Range("D2").Value = 5.36
$chor2.Range("I2").Value = 0.86
$chor2.Range("D3").Value = 5.77
$chor2.Range("I3").Value = 0.29
$chor2.Range("D4").Value = 4.2

# --- Add a new KONTROLA sheet (control group) with the same layout, placed after CHOR2 ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$kontrola = $wb.Worksheets.Add($null, $lastSheet)
$kontrola.Name = "KONTROLA"

$headers = @("group_name", "wiek", "hsCRP", "ERY", "PLT", "HGB", "HCT", "MCHC", "MON", "LEU")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $kontrola.Cells.Item(1, $i + 1).Value = $headers[$i]
}

$kontrola.Range("A2").Value = "KONTROLA sr"
$kontrola.Range("B2").Value = 29.56
$kontrola.Range("C2").Value = 6.1
$kontrola.Range("D2").Value = 5.36
$kontrola.Range("E2").Value = 225.28
$kontrola.Range("F2").Value = 12.41
$kontrola.Range("G2").Value = 0.36
$kontrola.Range("H2").Value = 35.13
$kontrola.Range("I2").Value = 0.86
$kontrola.Range("J2").Value = 12.02

$kontrola.Range("A3").Value = "KONTROLA os"
$kontrola.Range("B3").Value = 5.88
$kontrola.Range("C3").Value = 8.82
$kontrola.Range("D3").Value = 5.77
$kontrola.Range("E3").Value = 54.22
$kontrola.Range("F3").Value = 1.19
$kontrola.Range("G3").Value = 0.03
$kontrola.Range("H3").Value = 0.88
$kontrola.Range("I3").Value = 0.29
$kontrola.Range("J3").Value = 2.58

$kontrola.Range("A4").Value = "KONTROLA med"
$kontrola.Range("B4").Value = 29
$kontrola.Range("C4").Value = 3.97
$kontrola.Range("D4").Value = 4.2
$kontrola.Range("E4").Value = 217
$kontrola.Range("F4").Value = 12.4
$kontrola.Range("G4").Value = 0.36
$kontrola.Range("H4").Value = 35.05
$kontrola.Range("I4").Value = 0.76
$kontrola.Range("J4").Value = 11.66
